# Auto-generated edit script: update "想去人数" (F column) counts per commit 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1141   # was 1137
$ws.Range("F8").Value = 965   # was 953
$ws.Range("F10").Value = 6102   # was 6098
$ws.Range("F12").Value = 1760   # was 1754
$ws.Range("F13").Value = 450   # was 447
$ws.Range("F14").Value = 5998   # was 5979
$ws.Range("F15").Value = 117   # was 116
$ws.Range("F18").Value = 98   # was 97
$ws.Range("F19").Value = 1660   # was 1656
$ws.Range("F22").Value = 147   # was 146
$ws.Range("F23").Value = 1418   # was 1409
$ws.Range("F24").Value = 729   # was 726
$ws.Range("F25").Value = 255   # was 247

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 311   # was 310
$ws.Range("F5").Value = 170   # was 169
$ws.Range("F8").Value = 382   # was 381

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 9504   # was 9501
$ws.Range("F3").Value = 2235   # was 2234
$ws.Range("F5").Value = 196   # was 191

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 9504   # was 9501
$ws.Range("F3").Value = 2235   # was 2234
$ws.Range("F7").Value = 1141   # was 1137
$ws.Range("F11").Value = 311   # was 310
$ws.Range("F12").Value = 965   # was 953
$ws.Range("F13").Value = 196   # was 191
$ws.Range("F15").Value = 6102   # was 6098
$ws.Range("F17").Value = 1760   # was 1754
$ws.Range("F20").Value = 450   # was 447
$ws.Range("F23").Value = 5998   # was 5979
$ws.Range("F24").Value = 117   # was 116
$ws.Range("F27").Value = 98   # was 97
$ws.Range("F28").Value = 1660   # was 1656
$ws.Range("F31").Value = 147   # was 146
$ws.Range("F32").Value = 1418   # was 1409
$ws.Range("F33").Value = 729   # was 726
$ws.Range("F35").Value = 255   # was 247

